$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (old D:K shifts to F:M)
$ws.Columns("D:E").Insert()

# Copy number formatting (date / number styles) from the columns that used to be D:E (now F:G)
# onto the newly inserted D:E columns so they match the surrounding data (date row uses date
# format, everything else uses the plain number format), without creating duplicate style entries.
# Done in segments that skip the blank separator rows (36 and 78) so we don't materialize empty
# row/cell entries that didn't exist in the original sheet.
$ws.Range("F7:G35").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("F37:G77").Copy()
$ws.Range("D37").PasteSpecial(-4122)
$ws.Range("F79:G102").Copy()
$ws.Range("D79").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Approximate column widths for the two new columns using the widths of the old D/E (now F/G)
# columns.
$ws.Columns("D").ColumnWidth = $ws.Columns("F").ColumnWidth
$ws.Columns("E").ColumnWidth = $ws.Columns("G").ColumnWidth

# Populate the new columns with the latest two quarters of financial data.
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 29700
$ws.Range("E8").Value = 30500
$ws.Range("D9:E9").Value = "NA"
$ws.Range("D10:E10").Value = "NA"
$ws.Range("D12:E12").Value = "NA"
$ws.Range("D13:E13").Value = 0
$ws.Range("D14:E14").Value = 0
$ws.Range("D15:E15").Value = 0
$ws.Range("D17").Value = 200
$ws.Range("E17").Value = 300
$ws.Range("D18").Value = 29500
$ws.Range("E18").Value = 30200
$ws.Range("D20:E20").Value = 0
$ws.Range("D21:E21").Value = "NA"
$ws.Range("D22:E22").Value = 0
$ws.Range("D23").Value = 29500
$ws.Range("E23").Value = 30200
$ws.Range("D24:E24").Value = 0
$ws.Range("D25:E25").Value = 0
$ws.Range("D26").Value = 29500
$ws.Range("E26").Value = 30200
$ws.Range("D27").Value = 29500
$ws.Range("E27").Value = 30200
$ws.Range("D28:E28").Value = 0
$ws.Range("D29:E29").Value = 0
$ws.Range("D30:E30").Value = 0
$ws.Range("D31:E31").Value = 0
$ws.Range("D32:E32").Value = 0
$ws.Range("D33").Value = 29500
$ws.Range("E33").Value = 30200
$ws.Range("D34:E34").Value = 0
$ws.Range("D35").Value = 29500
$ws.Range("E35").Value = 30200
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41:E41").Value = 1000
$ws.Range("D42:E42").Value = 0
$ws.Range("D43:E43").Value = 0
$ws.Range("D44:E44").Value = 0
$ws.Range("D45:E45").Value = 0
$ws.Range("D46:E46").Value = 0
$ws.Range("D47:E47").Value = 0
$ws.Range("D48:E48").Value = 0
$ws.Range("D49:E49").Value = 0
$ws.Range("D50:E50").Value = 0
$ws.Range("D51:E51").Value = 0
$ws.Range("D52:E52").Value = 0
$ws.Range("D53:E53").Value = 0
$ws.Range("D54:E54").Value = 1000
$ws.Range("D57:E57").Value = 0
$ws.Range("D58:E58").Value = 0
$ws.Range("D59").Value = 300
$ws.Range("E59").Value = 200
$ws.Range("D60:E60").Value = 0
$ws.Range("D61:E61").Value = 0
$ws.Range("D62:E62").Value = 0
$ws.Range("D63:E63").Value = 0
$ws.Range("D64:E64").Value = 0
$ws.Range("D65:E65").Value = 0
$ws.Range("D66").Value = 300
$ws.Range("E66").Value = 200
$ws.Range("D68:E68").Value = 0
$ws.Range("D69:E69").Value = 0
$ws.Range("D70:E70").Value = 0
$ws.Range("D71:E71").Value = 0
$ws.Range("D72:E72").Value = 0
$ws.Range("D73:E73").Value = 0
$ws.Range("D74:E74").Value = 0
$ws.Range("D75:E75").Value = 0
$ws.Range("D76").Value = 700
$ws.Range("E76").Value = 800
$ws.Range("D77:E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 29500
$ws.Range("E81").Value = 30200
$ws.Range("D83:E83").Value = 0
$ws.Range("D84:E84").Value = 0
$ws.Range("D85:E85").Value = 0
$ws.Range("D86:E86").Value = 0
$ws.Range("D87:E87").Value = 0
$ws.Range("D88:E88").Value = 0
$ws.Range("D89").Value = 29400
$ws.Range("E89").Value = 30300
$ws.Range("D91:E91").Value = 0
$ws.Range("D92:E92").Value = 0
$ws.Range("D93:E93").Value = 0
$ws.Range("D94:E94").Value = 0
$ws.Range("D96").Value = -29500
$ws.Range("E96").Value = -30100
$ws.Range("D97:E97").Value = 0
$ws.Range("D98:E98").Value = 0
$ws.Range("D99:E99").Value = 0
$ws.Range("D100").Value = -29500
$ws.Range("E100").Value = -30100
$ws.Range("D101:E101").Value = 0
$ws.Range("D102").Value = -100
$ws.Range("E102").Value = 200
